$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: Collapse the three "CORE COMPETENCIES" detail paragraphs
# into a single condensed summary paragraph.
# ------------------------------------------------------------------
$bullet = [char]0x2022

# Find the heading paragraph and walk to the three paragraphs beneath it.
$coreHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "CORE COMPETENCIES") {
        $coreHeading = $p
        break
    }
}

$firstDetail = $coreHeading.Next()
$secondDetail = $firstDetail.Next()
$thirdDetail = $secondDetail.Next()

# Delete the 2nd and 3rd detail paragraphs entirely (bottom-up so the
# ranges of the earlier paragraphs remain valid while deleting).
$rngThird = $d.Range($thirdDetail.Range.Start, $thirdDetail.Range.End)
$rngThird.Delete()

$rngSecond = $d.Range($secondDetail.Range.Start, $secondDetail.Range.End)
$rngSecond.Delete()

# Replace the text of the remaining (first) detail paragraph with the
# condensed version.
$firstDetail.Range.Text = "Product Management & Strategy " + $bullet + " Technical Product Development " + $bullet + " Platform & Infrastructure"

# ------------------------------------------------------------------
# Change 2: Add a new "TECHNICAL SKILLS" section (heading + three
# paragraphs) right after the last bullet of "Technical Leadership &
# Management" and before the closing "For a more detailed..." line.
# ------------------------------------------------------------------
function Get-ClosingRange($doc) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13) -eq "For a more detailed, full description of my experience, please visit my LinkedIn and Personal Site.") {
            $r = $p.Range
            $r.Collapse(1)
            return $r
        }
    }
}

$newParagraphs = @(
    "TECHNICAL SKILLS",
    "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development",
    "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; API Development",
    "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Security & Compliance"
)

foreach ($txt in $newParagraphs) {
    $ip = Get-ClosingRange $d
    $ip.InsertBefore($txt + [char]13)
}

# Re-find the new heading paragraph and apply the Heading2 style to it.
$skillsHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "TECHNICAL SKILLS") {
        $skillsHeading = $p
        break
    }
}
$skillsHeading.Style = "Heading2"

Write-Host "Edit complete"
